# "tweaks after jan run" - clean up stray proofErr spell/gram-check markers
# (and the run-splits they forced) left behind by Word's spell/grammar
# checker, and fix "revise" -> "review".
#
# Word's Find/Replace merges the runs spanned by a match (and drops any
# w:proofErr markers that are now fully inside the replaced text), which is
# exactly the cleanup this commit performs, so we drive every edit through
# $d.Content.Find.Execute(...).

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $r = $d.Content
    $found = $r.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $found) {
        throw "Find failed for: $find"
    }
}

# 1. "Paste this content into the [Etherpad] at: " - drop spellStart/End around Etherpad
Replace-Text "Paste this content into the Etherpad at: " "Paste this content into the Etherpad at: "

# 2. "We[ do not really share ][data,][ ]we" - drop gramStart/End around "data,"
Replace-Text " do not really share data, " " do not really share data, "

# 3. "- I revise at least 4 articles a year:" -> "- I review at least 4 articles a year:"
Replace-Text "- I revise at least 4 articles a year:" "- I review at least 4 articles a year:"

# 4. " (e.g. journal club, carpentries, ReproducibiliTea)" - drop gramStart/End around
#    "e.g." and spellStart/End around "ReproducibiliTea"
Replace-Text " (e.g. journal club, carpentries, ReproducibiliTea)" " (e.g. journal club, carpentries, ReproducibiliTea)"

# 5. " to work[ ][any more]" - drop spellStart/End around "any more"
Replace-Text " any more" " any more"

# 6. "so you are [fairly confident][ ]they follow" - drop gramStart/End around "fairly confident"
Replace-Text "s data and notes, so you are fairly confident " "s data and notes, so you are fairly confident "

# 7 & 10. "How good was this [lesson][:]" (x2) - drop gramStart/End around "lesson:"
Replace-Text "How good was this lesson:" "How good was this lesson:"
Replace-Text "How good was this lesson:" "How good was this lesson:"

# 8. Titin/kDa paragraph - drop spellStart/End around "kDa"
Replace-Text "You need to do a western blot of the protein Titin, the largest protein in the body with a molecular weight of 3,800 kDa. You found a" "You need to do a western blot of the protein Titin, the largest protein in the body with a molecular weight of 3,800 kDa. You found a"

# 9. "Exercise 2- Public general record" merge into one run
Replace-Text "Exercise 2- Public general record" "Exercise 2- Public general record"
